$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "High Priority break-up" sheet so that the
#    untouched historical numbers survive under the new sheet name
#    "Major update - High Priority ". The copy is placed right after
#    the source sheet and becomes sheetId 6 / sheet6.xml.
# ---------------------------------------------------------------------
$orig = $wb.Worksheets.Item("High Priority break-up")
$orig.Copy($null, $orig)
$copy = $wb.Worksheets.Item($orig.Index + 1)
$copy.Name = "Major update - High Priority "

# ---------------------------------------------------------------------
# 2) Rename the original sheet and refresh it with the new interannual
#    breakup numbers (now 3 rows: header, "Trend New", "IUCN").
# ---------------------------------------------------------------------
$orig.Name = "Interannual update - High Pri"

$orig.Range("A2").Value = "Trend New"
$orig.Range("B2").Value = 59
$orig.Range("C2").Value = 57.3
$orig.Range("D2").Value = 59
$orig.Range("E2").Value = 78.7

$orig.Range("A3").Value = "IUCN"
$orig.Range("B3").Value = 44
$orig.Range("C3").Value = 42.7
$orig.Range("D3").Value = 16
$orig.Range("E3").Value = 21.3

# ---------------------------------------------------------------------
# 3) Update "Trends Status" numbers.
# ---------------------------------------------------------------------
$trends = $wb.Worksheets.Item("Trends Status")

$trends.Range("B2").Value = 0
$trends.Range("C2").Value = 1
$trends.Range("D2").Value = 0
$trends.Range("E2").Value = 3.2

$trends.Range("B3").Value = 1
$trends.Range("C3").Value = 5
$trends.Range("D3").Value = 10
$trends.Range("E3").Value = 16.1

$trends.Range("B4").Value = 3
$trends.Range("C4").Value = 21
$trends.Range("D4").Value = 30
$trends.Range("E4").Value = 67.7

$trends.Range("B5").Value = 5
$trends.Range("C5").Value = 2
$trends.Range("D5").Value = 50
$trends.Range("E5").Value = 6.5

$trends.Range("B6").Value = 1
$trends.Range("C6").Value = 2
$trends.Range("D6").Value = 10
$trends.Range("E6").Value = 6.5

$trends.Range("B7").Value = 32
$trends.Range("C7").Value = 105

$trends.Range("B8").Value = 524
$trends.Range("C8").Value = 430

# ---------------------------------------------------------------------
# 4) Update "Priority Status" numbers.
# ---------------------------------------------------------------------
$priority = $wb.Worksheets.Item("Priority Status")
$priority.Range("B2").Value = 103
$priority.Range("B3").Value = 286
$priority.Range("B4").Value = 554

# ---------------------------------------------------------------------
# 5) Update "Species qualification" numbers and label.
# ---------------------------------------------------------------------
$qual = $wb.Worksheets.Item("Species qualification")
$qual.Range("A2").Value = "SoIB Assessment"
$qual.Range("B2").Value = 566

$qual.Range("B3").Value = 42
$qual.Range("C3").Value = 10

$qual.Range("C4").Value = 31
